$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that differ between row 2 and row 3 and need to be swapped (plain
# numbers / plain text - safe to round-trip through .Value2 as-is).
$plainCols = @("A", "B", "D", "E", "F", "G", "H", "P", "Q", "R")

foreach ($col in $plainCols) {
    $addr2 = "$col" + "2"
    $addr3 = "$col" + "3"
    $val2 = $ws.Range($addr2).Value2
    $val3 = $ws.Range($addr3).Value2
    $ws.Range($addr2).Value2 = $val3
    $ws.Range($addr3).Value2 = $val2
}

# Columns that hold date-looking text (e.g. "2003-06-01") which must stay as
# literal text, not be auto-converted to a date serial number. Prefix with an
# apostrophe to force text, then strip the formatting change that the
# text-force leaves behind so the cell ends up identical to a plain text
# cell (no explicit style applied), matching the original layout.
$dateCols = @("Y", "AA")

foreach ($col in $dateCols) {
    $addr2 = "$col" + "2"
    $addr3 = "$col" + "3"
    $val2 = $ws.Range($addr2).Value2
    $val3 = $ws.Range($addr3).Value2
    $ws.Range($addr2).Value2 = "'" + $val3
    $ws.Range($addr3).Value2 = "'" + $val2
    $ws.Range($addr2).ClearFormats()
    $ws.Range($addr3).ClearFormats()
}
